$d = $word.ActiveDocument

# Locate the "LOB1004..." requisito paragraph - everything after it
# (the trailing blank paragraph, the "Ver no Jupiter..." line and the
# "© 2020 ..." copyright line) is being dropped by this edit.
$anchorStart = $d.Content.Duplicate
$foundStart = $anchorStart.Find.Execute(
    "LOB1004: Cálculo II (Requisito fraco)", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 0)

# Locate the final "© 2020 ..." copyright paragraph - the last paragraph
# that must be removed.
$anchorEnd = $d.Content.Duplicate
$foundEnd = $anchorEnd.Find.Execute(
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundStart -and $foundEnd) {
    # Step past the "LOB1004..." paragraph's own mark (start of the blank
    # paragraph that follows it) through the copyright paragraph's own
    # mark (so the copyright paragraph itself is removed too), leaving the
    # "LOB1004..." paragraph and everything after the old copyright
    # paragraph untouched.
    $deleteStart = $anchorStart.End + 1
    $deleteEnd = $anchorEnd.End + 1

    $killRange = $d.Range($deleteStart, $deleteEnd)
    $killRange.Delete()
}
